$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: D3's demand value is cleared out (was numeric 0)
$ws.Range("D3").Value = ""

# Row 7 gets relabeled from "Other" to "Biogas"; its value moves down to
# the new "Other" row below, so D7 becomes 0 here.
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 0

# Insert the new row 8 ("Other") that now carries the value which used
# to live on row 7, copying row 7's formatting for the label cell.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 84.7782351097898
